# Weekly price update: insert one new daily/weekly record as the most
# recent observation (row 230) for this series. All subsequent rows
# (previously 230-289) shift down by one (now 231-290).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 230..289 down to 231..290, leaving a fresh blank row at 230.
$ws.Rows("230:230").Insert()

# Populate the newly inserted row with the new observation. Every column
# except the date mirrors the (now shifted-down) row 231, which held the
# prior most-recent observation for this market/product combination.
$ws.Range("A230").Value2 = 9
$ws.Range("B230").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C230").Value2 = "Metropolitana"
$ws.Range("D230").Value2 = 44754
$ws.Range("E230").Value2 = 13
$ws.Range("F230").Value2 = 300000001
$ws.Range("G230").Value2 = "Rabanito"
$ws.Range("H230").Value2 = "Sin especificar"
$ws.Range("I230").Value2 = "Primera"
$ws.Range("J230").Value2 = 7000
$ws.Range("K230").Value2 = 2500
$ws.Range("L230").Value2 = 3000
$ws.Range("M230").Value2 = 2750
$ws.Range("N230").Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Range("O230").Value2 = "Provincia de Chacabuco"
$ws.Range("P230").Value2 = 28
$ws.Range("Q230").Value2 = 100
$ws.Range("R230").Value2 = "Hortaliza"
